{"js": "// Update the worksheet date title and every \"a\u00f7b=c, d\" answer cell in the\n// 5-column practice table. Cells are addressed by (row, col) in the table's\n// own grid (not by text search) because one new answer text happens to equal\n// another cell's *old* text (\"542\u00f77=77, 3\"): a sequential search-and-replace\n// would clobber the wrong cell once the first write landed. Writing through\n// `cell.getRange().insertText(text, \"Replace\")` (rather than\n// `cell.body.insertText`) replaces only the run's text and keeps the\n// existing paragraph/run formatting (alignment, font, size) untouched, which\n// matches the diff (only the <w:t> contents change).\nconst body = context.document.body;\n\n// Title paragraph: \"2024-07-16 Tuesday\" -> \"2024-07-17 Wednesday\".\nconst title = body.paragraphs.getFirst();\ntitle.getRange().insertText(\"2024-07-17 Wednesday\", \"Replace\");\nawait context.sync();\n\nconst table = body.tables.getFirst();\n\n// New answer values, keyed by (row, col) in the table's native grid \u2014\n// data lives in rows 0, 4, 8, 12, 16 (the rows in between are spacer rows).\nconst grid = [\n  [0, 0, \"137\u00f76=22, 5\"],\n  [0, 1, \"773\u00f79=85, 8\"],\n  [0, 2, \"774\u00f76=129, 0\"],\n  [0, 3, \"317\u00f72=158, 1\"],\n  [0, 4, \"526\u00f79=58, 4\"],\n\n  [4, 0, \"834\u00f77=119, 1\"],\n  [4, 1, \"710\u00f79=78, 8\"],\n  [4, 2, \"346\u00f79=38, 4\"],\n  [4, 3, \"878\u00f79=97, 5\"],\n  [4, 4, \"446\u00f77=63, 5\"],\n\n  [8, 0, \"810\u00f73=270, 0\"],\n  [8, 1, \"837\u00f78=104, 5\"],\n  [8, 2, \"245\u00f73=81, 2\"],\n  [8, 3, \"346\u00f74=86, 2\"],\n  [8, 4, \"464\u00f78=58, 0\"],\n\n  [12, 0, \"395\u00f79=43, 8\"],\n  [12, 1, \"834\u00f76=139, 0\"],\n  [12, 2, \"990\u00f79=110, 0\"],\n  [12, 3, \"542\u00f77=77, 3\"],\n  [12, 4, \"372\u00f78=46, 4\"],\n\n  [16, 0, \"416\u00f79=46, 2\"],\n  [16, 1, \"420\u00f74=105, 0\"],\n  [16, 2, \"681\u00f73=227, 0\"],\n  [16, 3, \"748\u00f74=187, 0\"],\n  [16, 4, \"318\u00f72=159, 0\"],\n];\n\nfor (const [row, col, newText] of grid) {\n  const cell = table.getCell(row, col);\n  cell.getRange().insertText(newText, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Update the worksheet date title and every \"a\u00f7b=c, d\" answer cell in the\n# 5-column practice table. Cells are addressed by their (row, col) position\n# in the table (Word COM is 1-indexed) rather than by text search, because\n# one new answer text happens to equal another cell's *old* text\n# (\"542\u00f77=77, 3\"): a sequential Find/Replace would clobber the wrong cell\n# once the first replacement landed. Writing straight to `Cell.Range.Text`\n# only rewrites the run text and leaves the existing paragraph/run\n# formatting (alignment, font, size) untouched, matching the diff (only the\n# <w:t> contents change).\n\n$d = $word.ActiveDocument\n\n# Title paragraph: \"2024-07-16 Tuesday\" -> \"2024-07-17 Wednesday\".\n$d.Paragraphs.Item(1).Range.Text = \"2024-07-17 Wednesday\"\n\n$tbl = $d.Tables.Item(1)\n\n# New answer values, keyed by (row, col) \u2014 1-indexed, matching Word COM.\n# Data lives in rows 1, 5, 9, 13, 17 (the rows in between are spacer rows).\n$tbl.Cell(1, 1).Range.Text = \"137\u00f76=22, 5\"\n$tbl.Cell(1, 2).Range.Text = \"773\u00f79=85, 8\"\n$tbl.Cell(1, 3).Range.Text = \"774\u00f76=129, 0\"\n$tbl.Cell(1, 4).Range.Text = \"317\u00f72=158, 1\"\n$tbl.Cell(1, 5).Range.Text = \"526\u00f79=58, 4\"\n\n$tbl.Cell(5, 1).Range.Text = \"834\u00f77=119, 1\"\n$tbl.Cell(5, 2).Range.Text = \"710\u00f79=78, 8\"\n$tbl.Cell(5, 3).Range.Text = \"346\u00f79=38, 4\"\n$tbl.Cell(5, 4).Range.Text = \"878\u00f79=97, 5\"\n$tbl.Cell(5, 5).Range.Text = \"446\u00f77=63, 5\"\n\n$tbl.Cell(9, 1).Range.Text = \"810\u00f73=270, 0\"\n$tbl.Cell(9, 2).Range.Text = \"837\u00f78=104, 5\"\n$tbl.Cell(9, 3).Range.Text = \"245\u00f73=81, 2\"\n$tbl.Cell(9, 4).Range.Text = \"346\u00f74=86, 2\"\n$tbl.Cell(9, 5).Range.Text = \"464\u00f78=58, 0\"\n\n$tbl.Cell(13, 1).Range.Text = \"395\u00f79=43, 8\"\n$tbl.Cell(13, 2).Range.Text = \"834\u00f76=139, 0\"\n$tbl.Cell(13, 3).Range.Text = \"990\u00f79=110, 0\"\n$tbl.Cell(13, 4).Range.Text = \"542\u00f77=77, 3\"\n$tbl.Cell(13, 5).Range.Text = \"372\u00f78=46, 4\"\n\n$tbl.Cell(17, 1).Range.Text = \"416\u00f79=46, 2\"\n$tbl.Cell(17, 2).Range.Text = \"420\u00f74=105, 0\"\n$tbl.Cell(17, 3).Range.Text = \"681\u00f73=227, 0\"\n$tbl.Cell(17, 4).Range.Text = \"748\u00f74=187, 0\"\n$tbl.Cell(17, 5).Range.Text = \"318\u00f72=159, 0\"\n"}
